$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 1: headers ---
$ws1.Range("A1").Value = "titre1"
$ws1.Range("B1").Value = "titre2"
$ws1.Range("C1").Value = "titre 3"

# --- Row 2 ---
$ws1.Range("A2").Value = "'2"
$ws1.Range("B2").Value = 3
$ws1.Range("C2").Value = "toto"

# --- Row 3 ---
$ws1.Range("A3").ClearContents()
$ws1.Range("B3").Value = 4
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("C3").Formula = "=3.7+B2"

# --- Selection on Feuil1 ---
$null = $ws1.Range("E2").Select()

# --- Add Feuil2 after Feuil1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Feuil2"
$null = $ws2.Range("H8").Select()

# Feuil2 is the active/selected tab in the final workbook.
$null = $ws2.Activate()
